$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 7394
$ws1.Range("F7").Value = 4760
$ws1.Range("F8").Value = 7005
$ws1.Range("F35").Value = 545
$ws1.Range("F43").Value = 21

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 32
$ws2.Range("F33").Value = 606

# Sheet "本地生活" (local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value = 1566
$ws3.Range("F9").Value = 2459

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 7394
$ws4.Range("F10").Value = 4760
$ws4.Range("F11").Value = 32
$ws4.Range("F12").Value = 0
$ws4.Range("F18").Value = 1566
$ws4.Range("F19").Value = 2459
$ws4.Range("F37").Value = 545
$ws4.Range("F38").Value = 606
$ws4.Range("F49").Value = 21
